# Case and Fatality Demographics Data Updated
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Cases by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B2").Value = 271
$ws.Range("B3").Value = 1326
$ws.Range("B4").Value = 3640
$ws.Range("B5").Value = 15492
$ws.Range("B6").Value = 17054
$ws.Range("B7").Value = 14955
$ws.Range("B8").Value = 12570
$ws.Range("B9").Value = 4525
$ws.Range("B10").Value = 3052
$ws.Range("B11").Value = 1845
$ws.Range("B12").Value = 1195
$ws.Range("B13").Value = 1881
$ws.Range("C21").Select()

# ---------------------------------------------------------------------------
# Sheet: Cases by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 26252
$ws.Range("B3").Value = 50670
$ws.Range("B4").Value = 897
$ws.Range("B2:B4").Select()

# ---------------------------------------------------------------------------
# Sheet: Cases by RaceEthnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 953
$ws.Range("B3").Value = 12924
$ws.Range("B4").Value = 28119
$ws.Range("B5").Value = 508
$ws.Range("B6").Value = 26773
$ws.Range("B7").Value = 8542
$ws.Range("E18").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B5").Value = 237
$ws.Range("B6").Value = 792
$ws.Range("B7").Value = 2340
$ws.Range("B8").Value = 5395
$ws.Range("B9").Value = 4526
$ws.Range("B10").Value = 5871
$ws.Range("B11").Value = 6505
$ws.Range("B12").Value = 6452
$ws.Range("B13").Value = 16339
$ws.Range("C23").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 20338
$ws.Range("B3").Value = 28169
$ws.Range("B15").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Race-Ethnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 1020
$ws.Range("B3").Value = 4845
$ws.Range("B4").Value = 22534
$ws.Range("B6").Value = 19822
$ws.Range("B13").Select()
